$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = "90.688.57"
$dCell.ClearFormats()
$ws.Range("E2").Value = "  -0.45%  "

$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = "3.145.36"
$dCell.ClearFormats()
$ws.Range("E3").Value = "  +1.17%  "

$dCell = $ws.Range("D4")
$dCell.NumberFormat = "@"
$dCell.Value = "0.999"
$dCell.ClearFormats()
$ws.Range("E4").Value = "  -0.14%  "

$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = "238.00"
$dCell.ClearFormats()
$ws.Range("E5").Value = "  +8.63%  "

$dCell = $ws.Range("D6")
$dCell.NumberFormat = "@"
$dCell.Value = "642.76"
$dCell.ClearFormats()
$ws.Range("E6").Value = "  +3.16%  "

$ws.Range("E7").Value = "  +11.09%  "

$dCell = $ws.Range("D8")
$dCell.NumberFormat = "@"
$dCell.Value = "0.365"
$dCell.ClearFormats()
$ws.Range("E8").Value = "  -3.70%  "

$ws.Range("E9").Value = "  -0.03%  "

$dCell = $ws.Range("D10")
$dCell.NumberFormat = "@"
$dCell.Value = "3.133.89"
$dCell.ClearFormats()
$ws.Range("E10").Value = "  +0.89%  "

$dCell = $ws.Range("D11")
$dCell.NumberFormat = "@"
$dCell.Value = "0.723"
$dCell.ClearFormats()
$ws.Range("E11").Value = "  -0.54%  "

$ws.Range("E12").Value = "  +3.03%  "

$dCell = $ws.Range("D13")
$dCell.NumberFormat = "@"
$dCell.Value = "36.59"
$dCell.ClearFormats()
$ws.Range("E13").Value = "  +6.17%  "

$dCell = $ws.Range("D14")
$dCell.NumberFormat = "@"
$dCell.Value = "0.0000247"
$dCell.ClearFormats()
$ws.Range("E14").Value = "  -3.22%  "

$ws.Range("E15").Value = "  +3.82%  "

$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"
$dCell.Value = "90.451.49"
$dCell.ClearFormats()
$ws.Range("E16").Value = "  -0.58%  "

$dCell = $ws.Range("D17")
$dCell.NumberFormat = "@"
$dCell.Value = "3.722.11"
$dCell.ClearFormats()
$ws.Range("E17").Value = "  +1.12%  "

$dCell = $ws.Range("D18")
$dCell.NumberFormat = "@"
$dCell.Value = "3.182.36"
$dCell.ClearFormats()
$ws.Range("E18").Value = "  +2.24%  "

$dCell = $ws.Range("D19")
$dCell.NumberFormat = "@"
$dCell.Value = "3.76"
$dCell.ClearFormats()
$ws.Range("E19").Value = "  +0.42%  "

$dCell = $ws.Range("D20")
$dCell.NumberFormat = "@"
$dCell.Value = "0.0000220"
$dCell.ClearFormats()
$ws.Range("E20").Value = "  -0.16%  "

$dCell = $ws.Range("D21")
$dCell.NumberFormat = "@"
$dCell.Value = "14.51"
$dCell.ClearFormats()
$ws.Range("E21").Value = "  +3.20%  "

$dCell = $ws.Range("D22")
$dCell.NumberFormat = "@"
$dCell.Value = "449.80"
$dCell.ClearFormats()
$ws.Range("E22").Value = "  +3.57%  "

$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value = "5.68"
$dCell.ClearFormats()
$ws.Range("E23").Value = "  +9.88%  "

$dCell = $ws.Range("D24")
$dCell.NumberFormat = "@"
$dCell.Value = "9.09"
$dCell.ClearFormats()
$ws.Range("E24").Value = "  +2.87%  "

$dCell = $ws.Range("D25")
$dCell.NumberFormat = "@"
$dCell.Value = "6.04"
$dCell.ClearFormats()
$ws.Range("E25").Value = "  -3.04%  "

$dCell = $ws.Range("D26")
$dCell.NumberFormat = "@"
$dCell.Value = "91.12"
$dCell.ClearFormats()
$ws.Range("E26").Value = "  +5.19%  "

$dCell = $ws.Range("D27")
$dCell.NumberFormat = "@"
$dCell.Value = "12.46"
$dCell.ClearFormats()
$ws.Range("E27").Value = "  +2.77%  "

$dCell = $ws.Range("D28")
$dCell.NumberFormat = "@"
$dCell.Value = "3.284.90"
$dCell.ClearFormats()
$ws.Range("E28").Value = "  +0.35%  "

$ws.Range("E29").Value = "  +0.06%  "

$dCell = $ws.Range("D30")
$dCell.NumberFormat = "@"
$dCell.Value = "9.76"
$dCell.ClearFormats()
$ws.Range("E30").Value = "  +7.38%  "

$dCell = $ws.Range("D31")
$dCell.NumberFormat = "@"
$dCell.Value = "0.160"
$dCell.ClearFormats()
$ws.Range("E31").Value = "  -4.01%  "

$ws.Range("E32").Value = "  +16.23%  "

$dCell = $ws.Range("D33")
$dCell.NumberFormat = "@"
$dCell.Value = "0.199"
$dCell.ClearFormats()
$ws.Range("E33").Value = "  +31.87%  "

$dCell = $ws.Range("D34")
$dCell.NumberFormat = "@"
$dCell.Value = "3.84"
$dCell.ClearFormats()
$ws.Range("E34").Value = "  +2.39%  "

$dCell = $ws.Range("D35")
$dCell.NumberFormat = "@"
$dCell.Value = "519.00"
$dCell.ClearFormats()
$ws.Range("E35").Value = "  -1.23%  "

$ws.Range("E36").Value = "  +3.98%  "

$ws.Range("E37").Value = "  +5.30%  "

$dCell = $ws.Range("D38")
$dCell.NumberFormat = "@"
$dCell.Value = "7.17"
$dCell.ClearFormats()
$ws.Range("E38").Value = "  +0.96%  "

$ws.Range("E39").Value = "  +1.89%  "

$dCell = $ws.Range("D40")
$dCell.NumberFormat = "@"
$dCell.Value = "0.810"
$dCell.ClearFormats()
$ws.Range("E40").Value = "  -10.23%  "

$dCell = $ws.Range("D41")
$dCell.NumberFormat = "@"
$dCell.Value = "0.424"
$dCell.ClearFormats()
$ws.Range("E41").Value = "  +5.97%  "

$ws.Range("E42").Value = "  -0.26%  "

$ws.Range("E43").Value = "  -1.32%  "

$dCell = $ws.Range("D45")
$dCell.NumberFormat = "@"
$dCell.Value = "3.26"
$dCell.ClearFormats()
$ws.Range("E45").Value = "  +42.49%  "

$dCell = $ws.Range("D47")
$dCell.NumberFormat = "@"
$dCell.Value = "0.706"
$dCell.ClearFormats()
$ws.Range("E47").Value = "  +13.06%  "

$dCell = $ws.Range("D48")
$dCell.NumberFormat = "@"
$dCell.Value = "150.19"
$dCell.ClearFormats()
$ws.Range("E48").Value = "  +1.74%  "

$dCell = $ws.Range("D49")
$dCell.NumberFormat = "@"
$dCell.Value = "46.01"
$dCell.ClearFormats()
$ws.Range("E49").Value = "  +4.38%  "

$dCell = $ws.Range("D50")
$dCell.NumberFormat = "@"
$dCell.Value = "4.63"
$dCell.ClearFormats()
$ws.Range("E50").Value = "  +9.50%  "

$ws.Range("E51").Value = "  +4.35%  "
